$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.517.24"
$ws.Range("E2").Value = "  -0.63%  "

$ws.Range("D3").Value = "2.949.47"
$ws.Range("E3").Value = "  -2.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.85%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  -0.29%  "

$ws.Range("D9").Value = "2.946.90"
$ws.Range("E9").Value = "  -1.98%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.19%  "

$ws.Range("E11").Value = "  -4.28%  "

$ws.Range("E12").Value = "  +1.21%  "

$ws.Range("E13").Value = "  -2.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.04%  "

$ws.Range("E15").Value = "  -0.48%  "

$ws.Range("D16").Value = "65.564.79"
$ws.Range("E16").Value = "  -0.58%  "

$ws.Range("D17").Value = "3.440.01"
$ws.Range("E17").Value = "  -1.99%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.51%  "

$ws.Range("D19").Value = "2.950.65"
$ws.Range("E19").Value = "  -2.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +13.43%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "446.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.695"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.06"
$ws.Range("D24").Style = "Normal"

$ws.Range("E25").Value = "  -3.29%  "

$ws.Range("E26").Value = "  -0.98%  "

$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.00%  "

$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.02%  "

$ws.Range("E29").Value = "  +7.19%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.67%  "

$ws.Range("E31").Value = "  -0.48%  "

$ws.Range("E32").Value = "  -0.94%  "

$ws.Range("E33").Value = "  +3.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.36%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.972"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.09%  "

$ws.Range("E37").Value = "  -1.69%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "45.86"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.84%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.09"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.36%  "

$ws.Range("E40").Value = "  -7.78%  "

$ws.Range("E41").Value = "  +0.21%  "

$ws.Range("E42").Value = "  -1.97%  "

$ws.Range("E43").Value = "  -6.87%  "

$ws.Range("E44").Value = "  +0.43%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "381.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.59%  "

$ws.Range("E46").Value = "  -1.56%  "

$ws.Range("D47").Value = "2.679.90"
$ws.Range("E47").Value = "  -4.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.38%  "

$ws.Range("E49").Value = "  +0.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.66%  "

$ws.Range("E51").Value = "  +1.47%  "
